# "update my history." -- append June's history rows for 6/7, 6/8 and 6/9
# (serials 40701-40703): Taxi + Taxi Tax, AM/PM, payback, onto the "June"
# sheet, rows 11-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("June")
$ws.Activate() | Out-Null

# Row 4 already has exactly the formatting pattern the new rows need:
# date-style (s=9) on column A and C, plain on B/D/E/F, nothing on G/H.
# Copy that formatting down across A:F for the 8 new rows first, so the
# new cells pick up style index 9 on A/C instead of creating new styles.
$ws.Range("A4:F4").Copy() | Out-Null
$ws.Range("A11:F18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$rows = @(
  @{r=11; a=40701; b="showa"; c="PM"; d="Taxi";     e="payback"; f=10},
  @{r=12; a=40701; b="showa"; c="PM"; d="Taxi Tax"; e="payback"; f=2},
  @{r=13; a=40702; b="showa"; c="AM"; d="Taxi";     e="payback"; f=75},
  @{r=14; a=40702; b="showa"; c="AM"; d="Taxi Tax"; e="payback"; f=2},
  @{r=15; a=40702; b="showa"; c="PM"; d="Taxi";     e="payback"; f=11},
  @{r=16; a=40702; b="showa"; c="PM"; d="Taxi Tax"; e="payback"; f=2},
  @{r=17; a=40703; b="showa"; c="AM"; d="Taxi";     e="payback"; f=64},
  @{r=18; a=40703; b="showa"; c="AM"; d="Taxi Tax"; e="payback"; f=2}
)

foreach ($row in $rows) {
  $r = $row.r
  $ws.Cells.Item($r, 1).Value = $row.a
  $ws.Cells.Item($r, 2).Value = $row.b
  $ws.Cells.Item($r, 3).Value = $row.c
  $ws.Cells.Item($r, 4).Value = $row.d
  $ws.Cells.Item($r, 5).Value = $row.e
  $ws.Cells.Item($r, 6).Value = $row.f
}

# Column F used to be a wide, manually sized 15.875 chars; the new data is
# short, so the author auto-fit it down to a narrow, best-fit width.
$ws.Columns.Item(6).ColumnWidth = 8.7

# Leave the cursor where the author's last action left it.
$ws.Range("G19").Select() | Out-Null

Write-Host "Added June rows 11-18 (6/7-6/9 taxi/taxi-tax payback entries)."
